$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93 (shifts existing rows 93-111 down to 94-112)
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record
$ws.Cells.Item(93, 1).Value = 1
$ws.Cells.Item(93, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(93, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(93, 4).Value = 44889
$ws.Cells.Item(93, 5).Value = 15
$ws.Cells.Item(93, 6).Value = 100112038
$ws.Cells.Item(93, 7).Value = "Cebollín baby"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 290
$ws.Cells.Item(93, 11).Value = 1700
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = 1803
$ws.Cells.Item(93, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(93, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(93, 16).Value = 902
$ws.Cells.Item(93, 17).Value = 2
$ws.Cells.Item(93, 18).Value = "Hortaliza"
